$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Num SNPs" (B), "Num Multiallelic" (C), "Num Invariants" (D) counts
# per chromosome row, reflecting new filtering parameters.
$data = @{
    2  = @(104238, 7332, 683245)
    3  = @(103048, 6937, 736819)
    4  = @(109888, 8372, 626633)
    5  = @(82161,  5778, 537914)
    6  = @(97546,  7571, 523135)
    7  = @(123996, 8191, 955530)
    8  = @(92488,  6829, 568958)
    9  = @(146636, 10611, 891010)
    10 = @(91529,  6354, 597264)
    11 = @(96497,  6577, 683599)
    12 = @(93300,  6134, 680939)
    13 = @(4191,   436,  19054)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# Update the saved selection to match the new state.
$ws.Range("I7").Select()
